$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '39.466.73'
$ws.Range("E2").Value = '  +1.75%  '

$ws.Range("D3").Value = '2.159.43'
$ws.Range("E3").Value = '  +2.77%  '

$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '227.84'
$ws.Range("E5").Value = '  -0.55%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.624'
$ws.Range("E6").Value = '  +0.90%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '64.11'
$ws.Range("E7").Value = '  +3.96%  '

$ws.Range("E8").Value = '  +0.04%  '

$ws.Range("E9").Value = '  +2.55%  '

$ws.Range("E10").Value = '  +1.32%  '

$ws.Range("E11").Value = '  -0.13%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.99'
$ws.Range("E12").Value = '  +3.42%  '

$ws.Range("D13").Value = '2.480.33'
$ws.Range("E13").Value = '  +2.92%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '22.17'
$ws.Range("E14").Value = '  +0.43%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.813'
$ws.Range("E15").Value = '  +0.43%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.54'
$ws.Range("E16").Value = '  +0.77%  '

$ws.Range("D17").Value = '2.159.58'
$ws.Range("E17").Value = '  -5.85%  '

$ws.Range("D18").Value = '39.413.85'
$ws.Range("E18").Value = '  +1.58%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '71.83'
$ws.Range("E19").Value = '  -0.23%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.12'
$ws.Range("E20").Value = '  +0.43%  '

$ws.Range("E21").Value = '  +1.26%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '231.24'
$ws.Range("E22").Value = '  +1.57%  '

$ws.Range("E23").Value = '  -0.02%  '

$ws.Range("B24").Value = 'PancakeSwap'
$ws.Range("C24").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.35'
$ws.Range("E24").Value = '  +0.53%  '

$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.29'
$ws.Range("E25").Value = '  -3.96%  '

$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '172.31'
$ws.Range("E26").Value = '  +0.34%  '

$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.51'
$ws.Range("E27").Value = '  -0.31%  '

$ws.Range("E28").Value = '  +1.05%  '

$ws.Range("E29").Value = '  +2.70%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.42'
$ws.Range("E30").Value = '  -0.48%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.67'
$ws.Range("E31").Value = '  +6.40%  '

$ws.Range("E32").Value = '  +0.67%  '

$ws.Range("E33").Value = '  +1.65%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.14'
$ws.Range("E34").Value = '  +9.53%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.75'
$ws.Range("E35").Value = '  -0.58%  '

$ws.Range("E36").Value = '  -0.70%  '

$ws.Range("E37").Value = '  +0.40%  '

$ws.Range("E38").Value = '  -0.45%  '

$ws.Range("E39").Value = '  +0.04%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '103.72'
$ws.Range("E40").Value = '  +2.14%  '

$ws.Range("E41").Value = '  +0.77%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '17.72'
$ws.Range("E42").Value = '  -2.46%  '

$ws.Range("D43").Value = '1.540.01'
$ws.Range("E43").Value = '  +0.49%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.19'
$ws.Range("E44").Value = '  +4.29%  '

$ws.Range("B45").Value = 'FTXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.35'
$ws.Range("E45").Value = '  +5.92%  '

$ws.Range("B46").Value = 'HuobiToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.83'
$ws.Range("E46").Value = '  +0.81%  '

$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0928'
$ws.Range("E47").Value = '  +2.02%  '

$ws.Range("E48").Value = '  +4.86%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.74'
$ws.Range("E49").Value = '  +0.02%  '

$ws.Range("D50").Value = '2.363.53'
$ws.Range("E50").Value = '  +3.08%  '

$ws.Range("E51").Value = '  -0.37%  '
